$wb = $excel.ActiveWorkbook

$sheetNames = @("WFA", "WFA (2)", "WFA (3)", "WFA (4)", "WFA (5)", "WFA (6)", "WFA (7)")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Insert a blank row above the table header (row 3), pushing everything down.
    $ws.Rows("3:3").Insert()

    # Re-anchor the table to its new location.
    $lo = $ws.ListObjects.Item(1)
    $lo.Resize($ws.Range("A4:I5"))
}
